# feat: add 2022-Q4 data
#
# The workbook contains a "总计" (totals) summary sheet plus one sheet per
# fiscal quarter (newest first). This change introduces a new quarter
# (2022-Q4) of data:
#   - a new row is added at the top of the quarterly history in "总计"
#   - a new quarterly detail sheet is appended for the oldest quarter that
#     now "falls off" the front of the list (2021-Q3), carrying the data
#     that used to live in the 2021-Q3 sheet
#   - every existing quarterly sheet's data cascades one slot older
#     (the sheet that used to show 2022-Q3 now shows what 2022-Q4 contains,
#     etc.) and the freed-up first quarterly sheet gets the brand new
#     2022-Q4 figures
#   - every quarterly sheet tab is renamed to the next quarter back

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Grab references to the existing quarterly sheets (by their current
#    names) before anything gets renamed.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$q3_22 = $wb.Worksheets.Item("2022-Q3")
$q2_22 = $wb.Worksheets.Item("2022-Q2")
$q1_22 = $wb.Worksheets.Item("2022-Q1")
$q4_21 = $wb.Worksheets.Item("2021-Q4")
$q3_21 = $wb.Worksheets.Item("2021-Q3")

# ---------------------------------------------------------------------
# 2. Append a brand new sheet after the last quarterly sheet; this will
#    become the new, oldest quarter (2021-Q3) and will receive the data
#    that currently lives in the 2021-Q3 sheet.
# ---------------------------------------------------------------------
$newOldest = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q3_21)

# ---------------------------------------------------------------------
# 3. Cascade the quarterly data one slot older, oldest first so nothing
#    gets overwritten before it has been copied onward.
# ---------------------------------------------------------------------
$q3_21.UsedRange.Copy($newOldest.Range("A1"))
$q4_21.UsedRange.Copy($q3_21.Range("A1"))
$q1_22.UsedRange.Copy($q4_21.Range("A1"))
$q2_22.UsedRange.Copy($q1_22.Range("A1"))
$q3_22.UsedRange.Copy($q2_22.Range("A1"))

# ---------------------------------------------------------------------
# 4. Put the brand new 2022-Q4 figures into the now-freed first
#    quarterly sheet (still referenced by $q3_22). The D/E/F/G columns
#    are stored as text in this workbook (like all the other quarterly
#    sheets), so a leading apostrophe keeps them as text instead of
#    being auto-converted to numbers.
# ---------------------------------------------------------------------
$q4sheet = $q3_22
$q4sheet.Cells.Item(2, 4).Value = "'4.76"
$q4sheet.Cells.Item(2, 5).Value = "'92.90"
$q4sheet.Cells.Item(2, 6).Value = "'5.20"
$q4sheet.Cells.Item(2, 7).Value = "'0.2475"
$q4sheet.Cells.Item(2, 8).Value = 5

$q4sheet.Cells.Item(3, 4).Value = "'0.64"
$q4sheet.Cells.Item(3, 5).Value = "'93.56"
$q4sheet.Cells.Item(3, 6).Value = "'3.80"
$q4sheet.Cells.Item(3, 7).Value = "'0.0243"

# ---------------------------------------------------------------------
# 5. Rename every quarterly tab to the quarter it now represents.
# ---------------------------------------------------------------------
$q4sheet.Name = "2022-Q4"
$q2_22.Name = "2022-Q3"
$q1_22.Name = "2022-Q2"
$q4_21.Name = "2022-Q1"
$q3_21.Name = "2021-Q4"
$newOldest.Name = "2021-Q3"

# ---------------------------------------------------------------------
# 6. Update the "总计" summary sheet: insert a new row for 2022-Q4 at the
#    top of the history table and renumber the index column.
# ---------------------------------------------------------------------
$summary.Rows.Item(2).Insert()

# The freshly inserted row inherits a blended format from its
# neighbours; clear B2:D2 back to the plain/default look used by the
# rest of the data rows, and give A2 the same style as the other index
# cells (A3:A7) by copying A3's formatting onto it.
$summary.Range("B2:D2").ClearFormats()
$summary.Cells.Item(3, 1).Copy($summary.Cells.Item(2, 1))

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q4"
$summary.Cells.Item(2, 3).Value = 2
$summary.Cells.Item(2, 4).Value = 0.27

$summary.Cells.Item(3, 1).Value = 1
$summary.Cells.Item(4, 1).Value = 2
$summary.Cells.Item(5, 1).Value = 3
$summary.Cells.Item(6, 1).Value = 4
$summary.Cells.Item(7, 1).Value = 5
